$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.343.30'
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").Value = '3.769.32'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'593.22"
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").Value = "'166.17"
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("D7").Value = '3.768.66'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("E11").Value = '  -1.46%  '
$ws.Range("D12").Value = "'0.450"
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("E13").Value = '  -3.13%  '
$ws.Range("D14").Value = "'35.86"
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("D15").Value = '4.407.72'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").Value = '3.783.05'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = '67.477.03'
$ws.Range("D18").Value = "'17.70"
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").Value = "'6.93"
$ws.Range("E20").Value = '  -1.72%  '
$ws.Range("E21").Value = '  -3.76%  '
$ws.Range("D22").Value = "'457.96"
$ws.Range("E22").Value = '  -2.15%  '
$ws.Range("D23").Value = "'0.695"
$ws.Range("E23").Value = '  -0.92%  '
$ws.Range("E24").Value = '  +7.08%  '
$ws.Range("D25").Value = "'83.24"
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("E26").Value = '  -4.52%  '
$ws.Range("E27").Value = '  -3.03%  '
$ws.Range("D28").Value = "'10.03"
$ws.Range("E28").Value = '  -1.74%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").Value = "'2.76"
$ws.Range("E30").Value = '  -1.33%  '
$ws.Range("D31").Value = "'7.20"
$ws.Range("E31").Value = '  -2.84%  '
$ws.Range("D32").Value = "'29.70"
$ws.Range("E32").Value = '  -1.13%  '
$ws.Range("D33").Value = "'2.16"
$ws.Range("E33").Value = '  -2.61%  '
$ws.Range("D34").Value = "'9.16"
$ws.Range("E34").Value = '  -1.36%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '3.722.71'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -1.61%  '
$ws.Range("E38").Value = '  -0.99%  '
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").Value = "'0.993"
$ws.Range("E40").Value = '  -1.01%  '
$ws.Range("D41").Value = "'5.73"
$ws.Range("E41").Value = '  -1.32%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = "'44.13"
$ws.Range("E44").Value = '  +1.57%  '
$ws.Range("D45").Value = "'0.300"
$ws.Range("E45").Value = '  -2.64%  '
$ws.Range("D46").Value = "'46.80"
$ws.Range("E46").Value = '  +3.01%  '
$ws.Range("E47").Value = '  -3.00%  '
$ws.Range("D48").Value = "'8.35"
$ws.Range("E48").Value = '  -2.58%  '
$ws.Range("D49").Value = "'146.10"
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").Value = "'390.92"
$ws.Range("E50").Value = '  -3.71%  '
$ws.Range("D51").Value = '2.753.81'
$ws.Range("E51").Value = '  +2.78%  '
